# Updated symbol list on Mon Dec 19 15:36:15 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are stored as text in this sheet, so force the
# number-format to Text before assigning each numeric-looking string --
# otherwise Excel/COM would auto-convert them to real numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "247.58"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.73"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.422"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05700"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8089"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.026"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1453"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07526"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03158"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03051"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09279"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.602"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001624"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04696"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006361"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005034"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001042"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0001500"
$ws.Range("E22").Value = "21UpBotsUBXT"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.772"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.409"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.099"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.3290"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006993"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002931"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008529"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005929"
$ws.Range("E47").Value = "46ACDXExchangeACXTWorstin24h"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.6827"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.007188"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
